$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the DC to DC converter part numbers
$ws.Range("B10").Value = "PYBE20-Q24-S5-T"
$ws.Range("B11").Value = "PYBE20-Q24-S24-T"

# Add note on the 24V converter row saying it's only needed for the Sonar
$ws.Range("D11").Value = "Only for the Sonar *2"

# Add footnote *2 explaining the sonar supply voltage compatibility
$ws.Range("A24").Value = "*2 you can verify with Imagenex if your sonar is comptatible with 12v supply"
$ws.Range("A25").Value = "if yes you can ommit the 24v suply"

# Fix Sonar component number typo: Imaginex 852 -> Imagenex 852
$ws.Range("B7").Value = "Imagenex 852"

# Restore the current selection to reflect the last-edited cell
$ws.Range("B7").Select()
